# Update Bill of Materials sheet to reflect newly added components (C26, C27, C28)
# and corrected resistor designator list / quantities.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-AndromedaV3.2")

# Row 2: CAP CER 10UF 10V X5R 0603 -> add C26, quantity 4 -> 5
# Leading "'" preserves the existing quotePrefix text-cell styling of column C.
$ws.Range("C2").Formula = "'C1, C17, C19, C20, C26"
$ws.Range("F2").Value = 5

# Row 3: CAP CER 0.1UF 25V X7R 0603 -> add C27, C28, quantity 14 -> 16
$ws.Range("C3").Formula = "'C2, C3, C4, C5, C12, C13, C14, C15, C16, C18, C21, C22, C23, C24, C27, C28"
$ws.Range("F3").Value = 16

# Row 19: RES SMD 10K OHM 5% 1/10W 0603 -> remove R9, R12, R13, quantity 13 -> 10
$ws.Range("C19").Formula = "'R2, R6, R7, R8, R11, R14, R15, R16, R17, R18"
$ws.Range("F19").Value = 10
